# Apply cryptocurrency list update (price/volume refresh + two coin-row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.225.40"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "2.253.71"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "274.05"
$ws.Range("E5").Value = "  +5.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.73"
$ws.Range("E6").Value = "  +10.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.46"
$ws.Range("E10").Value = "  +4.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  +8.70%  "

$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("D14").Value = "2.592.19"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.06"
$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").Value = "2.251.46"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.801"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "44.131.87"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("E20").Value = "  -0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.62"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.77"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").Value = "  -5.52%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.55"
$ws.Range("E26").Value = "  +14.25%  "

$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.56"
$ws.Range("E28").Value = "  +6.17%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.32"
$ws.Range("E29").Value = "  +5.36%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.66"
$ws.Range("E30").Value = "  -6.23%  "

$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.95"
$ws.Range("E32").Value = "  +1.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0899"
$ws.Range("E33").Value = "  +2.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").Value = "  +2.41%  "

$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("E36").Value = "  +2.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("E37").Value = "  -4.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.49"
$ws.Range("E39").Value = "  +21.97%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.22"
$ws.Range("E40").Value = "  +3.34%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.71"
$ws.Range("E41").Value = "  -4.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.02"
$ws.Range("E42").Value = "  +4.77%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.206"
$ws.Range("E43").Value = "  +0.49%  "

$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.47"
$ws.Range("E44").Value = "  +1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.53"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0992"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.08"
$ws.Range("E47").Value = "  -3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.21"
$ws.Range("E48").Value = "  +4.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.15"
$ws.Range("E49").Value = "  +1.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.433"
$ws.Range("E50").Value = "  -7.92%  "

$ws.Range("E51").Value = "  -0.03%  "
